# Applies the cryptos.xlsx price/volume refresh + 4 row-position swaps
# described in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so values such as "1.00", "7.01"
# or "69.293.51" are stored as literal strings (matching the workbook's
# original inline-string cells) instead of being auto-coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.293.51"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "3.682.37"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "682.47"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "158.63"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "0.147"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").Value = "7.01"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("D11").Value = "0.438"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("D13").Value = "4.304.74"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "32.40"
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").Value = "3.677.62"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "69.307.36"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").Value = "16.06"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("D19").Value = "6.43"
$ws.Range("E19").Value = "  -3.64%  "
$ws.Range("D20").Value = "471.40"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "10.01"
$ws.Range("E21").Value = "  +1.17%  "
$ws.Range("D22").Value = "0.652"
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "79.96"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "3.829.17"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  -6.79%  "
$ws.Range("D27").Value = "10.97"
$ws.Range("E27").Value = "  -4.71%  "
$ws.Range("D28").Value = "9.14"
$ws.Range("E28").Value = "  -5.91%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -5.93%  "
$ws.Range("D31").Value = "6.65"
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -5.44%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "26.92"
$ws.Range("E34").Value = "  -0.71%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.163"
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.664.79"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("D37").Value = "8.24"
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("D38").Value = "6.17"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D40").Value = "2.27"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "0.0906"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").Value = "174.12"
$ws.Range("E43").Value = "  +9.31%  "
$ws.Range("D44").Value = "0.943"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "47.56"
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").Value = "0.000286"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.72"
$ws.Range("E47").Value = "  -7.10%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "1.29"
$ws.Range("E48").Value = "  -5.22%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "27.64"
$ws.Range("E50").Value = "  -7.29%  "
$ws.Range("D51").Value = "7.80"
$ws.Range("E51").Value = "  -3.58%  "
